$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-21 Thursday" "2025-08-22 Friday"

Replace-Text "33×47=" "94×12="
Replace-Text "51×91=" "88×17="
Replace-Text "48×67=" "20×38="
Replace-Text "93×97=" "46×60="
Replace-Text "52×12=" "19×95="
Replace-Text "62×67=" "53×61="
Replace-Text "66×72=" "41×94="
Replace-Text "66×96=" "23×83="
Replace-Text "29×54=" "49×24="
Replace-Text "64×46=" "54×18="
Replace-Text "93×24=" "40×36="
Replace-Text "19×74=" "17×97="
Replace-Text "61×11=" "74×20="
Replace-Text "49×57=" "74×47="
Replace-Text "35×44=" "56×61="
Replace-Text "58×63=" "82×91="
Replace-Text "16×24=" "77×76="
Replace-Text "57×61=" "92×95="
Replace-Text "50×94=" "25×52="
Replace-Text "37×44=" "68×90="
Replace-Text "40×26=" "67×68="
Replace-Text "89×95=" "84×26="
Replace-Text "31×27=" "84×39="
Replace-Text "28×97=" "60×92="
Replace-Text "74×76=" "58×31="
